$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.962.69'
$ws.Range("E2").Value = '  +1.97%  '
$ws.Range("D3").Value = '1.814.01'
$ws.Range("E3").Value = '  +2.46%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.40'
$ws.Range("E5").Value = '  +2.21%  '
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4291'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3668'
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07247'
$ws.Range("E9").Value = '  +0.65%  '
$ws.Range("D10").Value = '2.196.33'
$ws.Range("E10").Value = '  +23.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8623'
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.25'
$ws.Range("E12").Value = '  +4.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.401'
$ws.Range("E13").Value = '  +3.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.602'
$ws.Range("E14").Value = '  +2.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06956'
$ws.Range("E15").Value = '  +0.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.12'
$ws.Range("E16").Value = '  +2.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.012'
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008892'
$ws.Range("E18").Value = '  +2.45%  '
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("E20").Value = '  +1.09%  '
$ws.Range("D21").Value = '27.012.15'
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.166'
$ws.Range("E22").Value = '  +1.27%  '
$ws.Range("D23").Value = '2.438.56'
$ws.Range("E23").Value = '  +22.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.01'
$ws.Range("E24").Value = '  -2.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.05'
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.865'
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.32'
$ws.Range("E27").Value = '  +1.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.228'
$ws.Range("E28").Value = '  +2.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.892'
$ws.Range("E29").Value = '  +8.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '114.44'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08925'
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.190'
$ws.Range("E32").Value = '  +6.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7461'
$ws.Range("E33").Value = '  +3.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.417'
$ws.Range("E34").Value = '  +2.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.806'
$ws.Range("E35").Value = '  +2.04%  '
$ws.Range("E36").Value = '  +0.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.124'
$ws.Range("E37").Value = '  +4.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05205'
$ws.Range("E38").Value = '  +0.99%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01921'
$ws.Range("E39").Value = '  +1.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5102'
$ws.Range("E40").Value = '  +3.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.756'
$ws.Range("E41").Value = '  +6.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1649'
$ws.Range("E42").Value = '  +2.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.477'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.320'
$ws.Range("E44").Value = '  +4.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '106.66'
$ws.Range("E45").Value = '  +1.74%  '
$ws.Range("E46").Value = '  +2.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.006'
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4564'
$ws.Range("E48").Value = '  +1.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.643'
$ws.Range("E49").Value = '  +3.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06212'
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.845'
$ws.Range("E51").Value = '  +6.05%  '
